$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row just under the data table: average of the "k" column (J).
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New labeled summary rows below that, with their computed statistics.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the four new statistic values: bold, 12pt, vertically centered.
$base = $ws.Cells.Item(14, 2)
$base.Font.Bold = $true
$base.Font.Size = 12
$base.VerticalAlignment = -4108

# Re-use that exact resolved style for the other three cells instead of
# re-deriving it property-by-property (avoids generating redundant style
# table entries).
$base.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Select the new block, matching the author's last on-screen selection.
$ws.Range("A14:B17").Select() | Out-Null

# Page setup tweak recorded alongside the edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
